$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 1.311598390708582
$arr[1,0] = 1.168688641629558
$arr[2,0] = 1.080820142514312
$arr[3,0] = 1.044984344992997
$arr[4,0] = 1.039032161265823
$arr[5,0] = 1.080336962042509
$arr[6,0] = 1.262349343736958
$arr[7,0] = 1.618250140937505
$arr[8,0] = 1.879049551238495
$arr[9,0] = 1.997536375466041
$arr[10,0] = 2.042381012937426
$arr[11,0] = 2.032724003744647
$arr[12,0] = 2.001226257417784
$arr[13,0] = 1.981929830090337
$arr[14,0] = 1.87130293631094
$arr[15,0] = 1.803396610921254
$arr[16,0] = 1.764324478939955
$arr[17,0] = 1.751092943506762
$arr[18,0] = 1.810626839512167
$arr[19,0] = 2.010478568418932
$arr[20,0] = 2.140953937465497
$arr[21,0] = 2.071330133617721
$arr[22,0] = 1.8073581523347
$arr[23,0] = 1.522085190698135
$ws.Range("B2:B25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 0.3005202967778189
$arr[1,0] = 0.262312094503045
$arr[2,0] = 0.2387608191005484
$arr[3,0] = 0.2291411008859541
$arr[4,0] = 0.2275424191884383
$arr[5,0] = 0.2386311740741291
$arr[6,0] = 0.2873653838335599
$arr[7,0] = 0.3821884631859689
$arr[8,0] = 0.4513810373882734
$arr[9,0] = 0.4827519515910694
$arr[10,0] = 0.4946157355812488
$arr[11,0] = 0.4920613648747576
$arr[12,0] = 0.483728310198785
$arr[13,0] = 0.4786220118327833
$arr[14,0] = 0.4493287088046145
$arr[15,0] = 0.4313308714013715
$arr[16,0] = 0.4209691297555196
$arr[17,0] = 0.4174591462287935
$arr[18,0] = 0.4332477968062562
$arr[19,0] = 0.4861763610482512
$arr[20,0] = 0.5206764602758085
$arr[21,0] = 0.5022716954971997
$arr[22,0] = 0.432381200099087
$arr[23,0] = 0.356618126746298
$ws.Range("C2:C25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 0.3346124636407239
$arr[1,0] = 0.3234465911793052
$arr[2,0] = 0.3166970551531136
$arr[3,0] = 0.313973437465819
$arr[4,0] = 0.3135228097385863
$arr[5,0] = 0.3166602145109465
$arr[6,0] = 0.3307404698984442
$arr[7,0] = 0.3591917247409526
$arr[8,0] = 0.3806041138337264
$arr[9,0] = 0.3904553423558923
$arr[10,0] = 0.3942015686634761
$arr[11,0] = 0.3933940520021793
$arr[12,0] = 0.39076323112252
$arr[13,0] = 0.3891538286117111
$arr[14,0] = 0.3799625279670806
$arr[15,0] = 0.3743522048962689
$arr[16,0] = 0.3711357155619055
$arr[17,0] = 0.3700484611321997
$arr[18,0] = 0.3749483558571285
$arr[19,0] = 0.3915355397358837
$arr[20,0] = 0.402468142412971
$arr[21,0] = 0.3966248377690818
$arr[22,0] = 0.3746788080900956
$arr[23,0] = 0.3514052827185594
$ws.Range("D2:D25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 1.190122002744133
$arr[1,0] = 1.193232468511795
$arr[2,0] = 1.196137627189934
$arr[3,0] = 1.197571135273634
$arr[4,0] = 1.197824226379367
$arr[5,0] = 1.196155950094919
$arr[6,0] = 1.190987519090967
$arr[7,0] = 1.188780480485704
$arr[8,0] = 1.19203713586829
$arr[9,0] = 1.194587731997487
$arr[10,0] = 1.1957080714743
$arr[11,0] = 1.195459903735269
$arr[12,0] = 1.194676802641709
$arr[13,0] = 1.194217271110304
$arr[14,0] = 1.191892025653104
$arr[15,0] = 1.190739921925399
$arr[16,0] = 1.190177849090645
$arr[17,0] = 1.190004794075008
$arr[18,0] = 1.190852149582483
$arr[19,0] = 1.194902619823495
$arr[20,0] = 1.198450633700233
$arr[21,0] = 1.19647432267881
$arr[22,0] = 1.190801099120463
$arr[23,0] = 1.18852413003232
$ws.Range("F2:F25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 0.00242772139124496
$arr[1,0] = 0.00243096625166983
$arr[2,0] = 0.002433063037825783
$arr[3,0] = 0.002433943835923612
$arr[4,0] = 0.002434091685296291
$arr[5,0] = 0.00243307480984067
$arr[6,0] = 0.002428818594244562
$arr[7,0] = 0.002421297031063738
$arr[8,0] = 0.002416268554526597
$arr[9,0] = 0.002414087904734604
$arr[10,0] = 0.002413277426679816
$arr[11,0] = 0.002413451298949089
$arr[12,0] = 0.002414020920348627
$arr[13,0] = 0.002414371818255173
$arr[14,0] = 0.002416413208321719
$arr[15,0] = 0.002417692842858651
$arr[16,0] = 0.002418438914727669
$arr[17,0] = 0.002418693251878007
$arr[18,0] = 0.002417555583106577
$arr[19,0] = 0.00241385319436178
$arr[20,0] = 0.002411522537861887
$arr[21,0] = 0.00241275832808942
$arr[22,0] = 0.002417617605982593
$arr[23,0] = 0.002423244049824302
$ws.Range("G2:G25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 0.3481897139701431
$arr[1,0] = 0.3367139205240051
$arr[2,0] = 0.3298786986282209
$arr[3,0] = 0.3271462053172769
$arr[4,0] = 0.3266956689824809
$arr[5,0] = 0.32984163319
$arr[6,0] = 0.3441889906122242
$arr[7,0] = 0.3740059142762249
$arr[8,0] = 0.3969513992574747
$arr[9,0] = 0.4076186414031469
$arr[10,0] = 0.4116912082536146
$arr[11,0] = 0.4108126329469997
$arr[12,0] = 0.4079530294072811
$arr[13,0] = 0.4062057565757158
$arr[14,0] = 0.3962588954138226
$arr[15,0] = 0.3902156357722504
$arr[16,0] = 0.3867612793695372
$arr[17,0] = 0.3855953931570468
$arr[18,0] = 0.3908567179156108
$arr[19,0] = 0.4087920643737419
$arr[20,0] = 0.4207069336947455
$arr[21,0] = 0.4143300261512479
$arr[22,0] = 0.3905668224155363
$arr[23,0] = 0.3657581849121527
$ws.Range("J2:J25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 2.474416872065831
$arr[1,0] = 2.500455950125939
$arr[2,0] = 2.518845044920113
$arr[3,0] = 2.526940535921796
$arr[4,0] = 2.528321074830217
$arr[5,0] = 2.518951789554947
$arr[6,0] = 2.482895636946296
$arr[7,0] = 2.4313300778837
$arr[8,0] = 2.405238934135951
$arr[9,0] = 2.395957484163432
$arr[10,0] = 2.392817111461227
$arr[11,0] = 2.393476763890959
$arr[12,0] = 2.395691609315577
$arr[13,0] = 2.397097078473053
$arr[14,0] = 2.405897756937719
$arr[15,0] = 2.411961003431259
$arr[16,0] = 2.415691841110771
$arr[17,0] = 2.416996778234534
$arr[18,0] = 2.41129035143635
$arr[19,0] = 2.395030879097249
$arr[20,0] = 2.386587177516077
$arr[21,0] = 2.390893278012015
$arr[22,0] = 2.411592790201325
$arr[23,0] = 2.443217832025681
$ws.Range("O2:O25").Value = $arr

